$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.250.17"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.861.31"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'0.7019"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'237.70"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.08110"
$ws.Range("E8").Value = "  +8.42%  "
$ws.Range("D9").Value = "'0.3026"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'23.19"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.08157"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.844.15"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "'5.158"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "'0.7049"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "'88.98"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "29.255.97"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'5.763"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'0.000007832"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "'13.30"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").Value = "'234.97"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D22").Value = "2.109.87"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'7.393"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").Value = "'161.25"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "'8.948"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'0.1441"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'18.05"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'1.962"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "'1.433"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "'1.484"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "'4.385"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'4.047"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'1.166"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").Value = "'0.7051"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").Value = "'0.9972"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").Value = "'2.676"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'0.01838"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").Value = "'2.730"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").Value = "'0.9203"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").Value = "1.131.81"
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("D43").Value = "'0.4264"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'5.867"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D45").Value = "'70.08"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'102.18"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "'1.761"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "2.002.86"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").Value = "'9.158"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "'6.931"
$ws.Range("E51").Value = "  -1.02%  "
